# Update the "3T" sheet (3rd quarter grades) data table, rows 2-29, columns B:L
$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

$sheet3Data = @(
    @(2, 5.3, 1.8, 4.9, 8.4, 2.2, 1.6, 1.5, 0, 5.6, 0.7, 3.2),
    @(3, 8.5, 5.1, 7.8, 2.9, 3.5, 4.6, 6.8, 8.6, 6.4, 0, 3.6),
    @(4, 5.8, 8.6, 7.3, 5.2, 9.2, 6.6, 6.8, 2.2, 7.8, 1.5, 6),
    @(5, 6.6, 8.2, 5.6, 6.8, 5.4, 8.2, 8.1, 9.8, 5, 2.5, 7.2),
    @(6, 8.6, 8.8, 4, 5.4, 6.2, 6.6, 7.8, 6, 6.6, 3.8, 3),
    @(7, 3.5, 4.9, 4, 7.5, 4.6, 1.6, 2.2, 1.2, 1.6, 0, 5.9),
    @(8, 3.4, 4.4, 6.4, 1.6, 0.1, 0.4, 6.4, 5.4, 6.7, 0, 6.8),
    @(9, 6.4, 8.3, 2.2, 8.2, 4.5, 6.2, 9.6, 2.3, 6.1, 4.5, 8),
    @(10, 7.4, 4, 2.7, 5.2, 4.2, 8.2, 6.5, 8.6, 2.6, 2.2, 4.4),
    @(11, 9.9, 7.6, 6.4, 8.1, 4.2, 2.9, 7.1, 8.8, 7.1, 3.7, 5.8),
    @(12, 7.4, 7.3, 3, 4.6, 2.4, 7.2, 4, 6.7, 7.4, 0, 9.2),
    @(13, 6, 9.2, 9.8, 6.6, 4.6, 4.6, 1.2, 4.6, 4.6, 3.8, 4.4),
    @(14, 9.3, 9.8, 7.2, 7.7, 8.6, 1.6, 7.8, 3.6, 5.8, 2.1, 3.4),
    @(15, 7.8, 5.4, 7.8, 6.4, 9.5, 8.4, 5.8, 8.2, 5.4, 5.8, 9.9),
    @(16, 4.8, 7.8, 7.4, 7.4, 8.7, 8.8, 9.9, 5.7, 8.9, 6.9, 4.5),
    @(17, 6.8, 3.1, 8.8, 8, 9.6, 8.1, 7.2, 5.7, 8.5, 0.4, 7.2),
    @(18, 9.6, 9.1, 6.1, 5.7, 4.4, 8.8, 2.8, 5, 6.3, 0.6, 6.4),
    @(19, 6.9, 5.8, 2.7, 1.5, 4.3, 1.1, 4.1, 4.6, 8.7, 0.7, 6.2),
    @(20, 7.1, 10, 5.2, 8, 10, 8.2, 9, 4.1, 8.6, 4.4, 8.4),
    @(21, 6.8, 4.9, 4.8, 0.6, 4.8, 8.4, 7.4, 2.2, 5, 2.6, 1),
    @(22, 8.7, 8.5, 4.6, 7.2, 5.6, 6.4, 2.9, 7.1, 8.8, 6, 9),
    @(23, 7.4, 8.2, 4.2, 3.3, 7.6, 0.9, 8.2, 6.1, 8.4, 3.2, 7),
    @(24, 6, 7.2, 8.8, 4.5, 8.4, 7.8, 5, 6.9, 9, 0, 7.2),
    @(25, 10, 8.2, 7.1, 5.7, 7.2, 3.4, 7, 2.6, 5.6, 3.4, 9.9),
    @(26, 8.5, 5.8, 5.8, 5.6, 9.9, 6.1, 8.8, 6.4, 8.8, 1.3, 6.6),
    @(27, 9.9, 9, 8.2, 9, 2.4, 8.2, 7, 6.4, 7, 0.2, 5.9),
    @(28, 3.7, 3.7, 9.1, 4.9, 8.5, 3.6, 6.4, 1, 6.4, 7.8, 4.8),
    @(29, 7.4, 8.4, 5.8, 8.8, 6, 7.8, 4, 7.2, 7.8, 4.4, 5.6)
)

foreach ($rowData in $sheet3Data) {
    $r = $rowData[0]
    $arr = New-Object 'object[,]' 1,11
    for ($i = 0; $i -lt 11; $i++) {
        $arr[0, $i] = $rowData[$i + 1]
    }
    $ws3.Range("B$r`:L$r").Value = $arr
}

# Update the "Média por Trimestre" (averages) sheet:
# 1) Add the missing "Trimestres" header label in A1, copying the existing
#    header formatting (bold, centered, bordered) from B1.
# 2) Refresh the 3º Trimestre row (row 4) with the new column averages.
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws4.Range("A1").Value = "Trimestres"

$row4Vals = @(7.124999999999999, 6.896428571428571, 5.989285714285715, 5.885714285714286, 5.950000000000001, 5.582142857142857, 6.117857142857145, 5.25, 6.660714285714286, 2.589285714285715, 6.089285714285715)
$row4Arr = New-Object 'object[,]' 1,11
for ($i = 0; $i -lt 11; $i++) {
    $row4Arr[0, $i] = $row4Vals[$i]
}
$ws4.Range("B4:L4").Value = $row4Arr
